$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A4 value from 1 to 10
$ws.Range("A4").Value = 10

# Add new row 5
$ws.Range("A5").Value = 20
$ws.Range("B5").Value = "first 20 in data folder"
$ws.Range("C5").Value = 96.2

# Update selection to B7
$ws.Range("B7").Select()
